$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "EyesClosed"
$ws.Range("B3").Value = "Attention"
$ws.Range("B4").Value = "Attention"
$ws.Range("B5").Value = "Attention"
$ws.Range("B6").Value = "Attention"
$ws.Range("B8").Value = "Attention"
$ws.Range("B9").Value = "EyesClosed"
$ws.Range("B10").Value = "Attention"
$ws.Range("B12").Value = "EyesClosed"
$ws.Range("B13").Value = "EyesClosed"
$ws.Range("B17").Value = "Attention"
$ws.Range("B18").Value = "Attention"
$ws.Range("B22").Value = "Attention"
$ws.Range("B25").Value = "Attention"
$ws.Range("B26").Value = "Attention"
$ws.Range("B27").Value = "Attention"
$ws.Range("B28").Value = "Attention"
$ws.Range("B30").Value = "EyesClosed"
$ws.Range("B33").Value = "EyesClosed"
$ws.Range("B34").Value = "EyesClosed"
$ws.Range("B39").Value = "EyesClosed"
$ws.Range("B45").Value = "EyesClosed"
$ws.Range("B47").Value = "Attention"
$ws.Range("B48").Value = "EyesClosed"
$ws.Range("B49").Value = "EyesClosed"
$ws.Range("B53").Value = "Attention"
$ws.Range("B54").Value = "EyesClosed"
$ws.Range("B57").Value = "Attention"
$ws.Range("B59").Value = "Attention"
$ws.Range("B60").Value = "Attention"
$ws.Range("B62").Value = "Attention"
$ws.Range("B65").Value = "EyesClosed"
$ws.Range("B67").Value = "EyesClosed"
$ws.Range("B68").Value = "EyesClosed"
$ws.Range("B69").Value = "Attention"
$ws.Range("B72").Value = "EyesClosed"
$ws.Range("B73").Value = "EyesClosed"
$ws.Range("B79").Value = "Attention"
$ws.Range("B82").Value = "Attention"
$ws.Range("B83").Value = "Attention"
$ws.Range("B84").Value = "Attention"
$ws.Range("B94").Value = "EyesClosed"
$ws.Range("B95").Value = "EyesClosed"
$ws.Range("B99").Value = "Attention"
$ws.Range("B100").Value = "Attention"
$ws.Range("B101").Value = "Attention"
$ws.Range("B102").Value = "Attention"
$ws.Range("B103").Value = "Attention"
$ws.Range("B104").Value = "Attention"
$ws.Range("B105").Value = "Attention"
$ws.Range("B106").Value = "Attention"
$ws.Range("B107").Value = "Attention"
$ws.Range("B108").Value = "Attention"
$ws.Range("B109").Value = "Attention"
$ws.Range("B110").Value = "Attention"
$ws.Range("B111").Value = "Attention"
$ws.Range("B112").Value = "Attention"
$ws.Range("B124").Value = "Attention"
$ws.Range("B137").Value = "Attention"
$ws.Range("B141").Value = "EyesClosed"
$ws.Range("B142").Value = "Attention"
$ws.Range("B144").Value = "Attention"
$ws.Range("B145").Value = "Attention"
$ws.Range("B152").Value = "Attention"
$ws.Range("B155").Value = "Attention"
$ws.Range("B157").Value = "Attention"
$ws.Range("B163").Value = "Attention"
$ws.Range("B169").Value = "Attention"
$ws.Range("B170").Value = "Attention"
$ws.Range("B171").Value = "EyesClosed"
$ws.Range("B175").Value = "EyesClosed"
$ws.Range("B176").Value = "EyesClosed"
$ws.Range("B181").Value = "EyesClosed"
$ws.Range("B183").Value = "EyesClosed"
$ws.Range("B185").Value = "Attention"
$ws.Range("B190").Value = "Attention"
$ws.Range("B191").Value = "Attention"
$ws.Range("B192").Value = "Attention"
$ws.Range("B193").Value = "Attention"
$ws.Range("B194").Value = "Attention"
$ws.Range("B195").Value = "Attention"
$ws.Range("B196").Value = "Attention"
$ws.Range("B197").Value = "Attention"
$ws.Range("B198").Value = "Attention"
$ws.Range("B199").Value = "Attention"
$ws.Range("B200").Value = "Attention"
$ws.Range("B204").Value = "Attention"
$ws.Range("B206").Value = "Attention"
$ws.Range("B211").Value = "Attention"
$ws.Range("B222").Value = "EyesClosed"
$ws.Range("B223").Value = "EyesClosed"
$ws.Range("B224").Value = "EyesClosed"
$ws.Range("B225").Value = "EyesClosed"
$ws.Range("B226").Value = "EyesClosed"
$ws.Range("B228").Value = "Attention"
$ws.Range("B229").Value = "Attention"
$ws.Range("B230").Value = "Attention"
$ws.Range("B231").Value = "Attention"
$ws.Range("B232").Value = "Attention"
$ws.Range("B233").Value = "Attention"
$ws.Range("B238").Value = "Attention"
$ws.Range("B239").Value = "Attention"
$ws.Range("B241").Value = "Attention"
